# Auto-generated: apply scheduled market-data refresh to Phoenix_Profits workbook
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 1958.3334
$ws.Range("I43").Value = 1917.8572
$ws.Range("K43").Value = 1917.8572
$ws.Range("M43").Value = -1848.8572
$ws.Range("H88").Value = 1284.1666
$ws.Range("I88").Value = 1152.6666
$ws.Range("J88").Value = 1328
$ws.Range("K88").Value = 1152.6666
$ws.Range("L88").Value = 1328
$ws.Range("M88").Value = -746.6666
$ws.Range("N88").Value = -2140
$ws.Range("H91").Value = 1284.1666
$ws.Range("I91").Value = 1152.6666
$ws.Range("J91").Value = 1328
$ws.Range("K91").Value = 1152.6666
$ws.Range("L91").Value = 1328
$ws.Range("M91").Value = 251.3334
$ws.Range("N91").Value = -4136
$ws.Range("H113").Value = 2384.3076
$ws.Range("I113").Value = 1855.2222
$ws.Range("J113").Value = 3574.75
$ws.Range("K113").Value = 1855.2222
$ws.Range("L113").Value = 3574.75
$ws.Range("M113").Value = 1398.7778
$ws.Range("N113").Value = -10082.75
$ws.Range("H116").Value = 7294.436
$ws.Range("I116").Value = 7000.125
$ws.Range("K116").Value = 7000.125
$ws.Range("M116").Value = -3558.125
$ws.Range("H137").Value = 2110.8572
$ws.Range("I137").Value = 2008.8334
$ws.Range("K137").Value = 6026.5002
$ws.Range("M137").Value = -3476.5002
$ws.Range("H138").Value = 2506.3674
$ws.Range("I138").Value = 1841.7084
$ws.Range("J138").Value = 3144.44
$ws.Range("K138").Value = 5525.1252
$ws.Range("L138").Value = 9433.32
$ws.Range("M138").Value = -385.1252000000004
$ws.Range("N138").Value = -19713.32

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2634.4348
$ws.Range("I2").Value = 2494.4736
$ws.Range("K2").Value = 2494.4736
$ws.Range("M2").Value = -2381.4736
$ws.Range("H32").Value = 2312.7974
$ws.Range("I32").Value = 2186.5078
$ws.Range("K32").Value = 2186.5078
$ws.Range("M32").Value = -1899.5078
$ws.Range("H34").Value = 40124.875
$ws.Range("I34").Value = 36666.332
$ws.Range("J34").Value = 42200
$ws.Range("K34").Value = 36666.332
$ws.Range("L34").Value = 42200
$ws.Range("M34").Value = -36395.332
$ws.Range("N34").Value = -42742
$ws.Range("H45").Value = 2465.7083
$ws.Range("I45").Value = 2062.0667
$ws.Range("K45").Value = 2062.0667
$ws.Range("M45").Value = -1685.0667
$ws.Range("H74").Value = 1905.5
$ws.Range("I74").Value = 1600
$ws.Range("K74").Value = 1600
$ws.Range("M74").Value = -726
$ws.Range("H77").Value = 1905.5
$ws.Range("I77").Value = 1600
$ws.Range("K77").Value = 8000
$ws.Range("M77").Value = -3632
$ws.Range("H80").Value = 41740.5
$ws.Range("J80").Value = 47898.25
$ws.Range("L80").Value = 47898.25
$ws.Range("N80").Value = -49894.25
$ws.Range("H83").Value = 41740.5
$ws.Range("J83").Value = 47898.25
$ws.Range("L83").Value = 143694.75
$ws.Range("N83").Value = -153678.75
$ws.Range("H116").Value = 2634.4348
$ws.Range("I116").Value = 2494.4736
$ws.Range("K116").Value = 2494.4736
$ws.Range("M116").Value = -200.4735999999998
$ws.Range("H132").Value = 10050
$ws.Range("I132").Value = 10050
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 30150
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -27620

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2634.4348
$ws.Range("I3").Value = 2494.4736
$ws.Range("K3").Value = 2494.4736
$ws.Range("M3").Value = -2380.4736
$ws.Range("H86").Value = 2434.5
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 2434.5
$ws.Range("K86").Value = 0
$ws.Range("L86").ClearContents()
$ws.Range("M86").Value = 2434.5
$ws.Range("N86").Value = -4680.5
$ws.Range("H89").Value = 2434.5
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 2434.5
$ws.Range("K89").Value = 0
$ws.Range("L89").ClearContents()
$ws.Range("M89").Value = 12172.5
$ws.Range("N89").Value = -23404.5
$ws.Range("H105").Value = 33337228
$ws.Range("I105").Value = 55559440
$ws.Range("K105").Value = 55559440
$ws.Range("M105").Value = -55557693

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2235.1606
$ws.Range("I31").Value = 1035.762
$ws.Range("K31").Value = 1035.762
$ws.Range("M31").Value = -740.7619999999999
$ws.Range("H34").Value = 2235.1606
$ws.Range("I34").Value = 1035.762
$ws.Range("K34").Value = 1035.762
$ws.Range("M34").Value = -833.7619999999999
$ws.Range("H122").Value = 7112.5293
$ws.Range("I122").Value = 7170.4614
$ws.Range("K122").Value = 21511.3842
$ws.Range("M122").Value = -19061.3842
$ws.Range("H125").Value = 49999.5
$ws.Range("J125").Value = 49999.5
$ws.Range("L125").Value = 49999.5
$ws.Range("N125").Value = -54919.5
$ws.Range("H133").Value = 93041
$ws.Range("I133").Value = 89948.5
$ws.Range("K133").Value = 89948.5
$ws.Range("M133").Value = -87418.5
$ws.Range("H134").Value = 2083.606
$ws.Range("I134").Value = 1891.5
$ws.Range("J134").Value = 4004.6667
$ws.Range("K134").Value = 5674.5
$ws.Range("L134").Value = 12014.0001
$ws.Range("M134").Value = -3139.5
$ws.Range("N134").Value = -17084.0001

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 1179.5
$ws.Range("I12").Value = 142.25
$ws.Range("J12").Value = 1368.091
$ws.Range("K12").Value = 426.75
$ws.Range("L12").Value = 4104.272999999999
$ws.Range("M12").Value = -253.75
$ws.Range("N12").Value = -4450.272999999999
$ws.Range("H107").Value = 1569.3438
$ws.Range("J107").Value = 1720.875
$ws.Range("L107").Value = 5162.625
$ws.Range("N107").Value = -9002.625
$ws.Range("H123").Value = 3999
$ws.Range("I123").Value = 3999
$ws.Range("K123").Value = 11997
$ws.Range("M123").Value = -9547
$ws.Range("H132").Value = 1658.619
$ws.Range("J132").Value = 1834.2307
$ws.Range("L132").Value = 16508.0763
$ws.Range("N132").Value = -21568.0763

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H33").Value = 0
$ws.Range("J33").Value = 0
$ws.Range("L33").ClearContents()
$ws.Range("N33").Value = 0
$ws.Range("H80").Value = 3142.2856
$ws.Range("I80").Value = 3329.6667
$ws.Range("J80").Value = 3001.75
$ws.Range("K80").Value = 3329.6667
$ws.Range("L80").Value = 3001.75
$ws.Range("M80").Value = -2331.6667
$ws.Range("N80").Value = -4997.75
$ws.Range("H83").Value = 3142.2856
$ws.Range("I83").Value = 3329.6667
$ws.Range("J83").Value = 3001.75
$ws.Range("K83").Value = 16648.3335
$ws.Range("L83").Value = 15008.75
$ws.Range("M83").Value = -11656.3335
$ws.Range("N83").Value = -24992.75
$ws.Range("H97").Value = 1249.3871
$ws.Range("I97").Value = 1342.4736
$ws.Range("J97").Value = 1102
$ws.Range("K97").Value = 1342.4736
$ws.Range("L97").Value = 1102
$ws.Range("M97").Value = -846.4736
$ws.Range("N97").Value = -2094
$ws.Range("H102").Value = 4354.636
$ws.Range("I102").Value = 3775.3572
$ws.Range("K102").Value = 3775.3572
$ws.Range("M102").Value = -2153.3572
$ws.Range("H126").Value = 3127.0667
$ws.Range("J126").Value = 2824.25
$ws.Range("L126").Value = 8472.75
$ws.Range("N126").Value = -13412.75
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("L132").ClearContents()
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = 0

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1298.2
$ws.Range("I16").Value = 1487.7333
$ws.Range("K16").Value = 1487.7333
$ws.Range("M16").Value = -1317.7333
$ws.Range("H68").Value = 2929.8125
$ws.Range("I68").Value = 1907.909
$ws.Range("K68").Value = 1907.909
$ws.Range("M68").Value = -1158.909
$ws.Range("H71").Value = 2929.8125
$ws.Range("I71").Value = 1907.909
$ws.Range("K71").Value = 9539.545
$ws.Range("M71").Value = -5795.545
$ws.Range("H82").Value = 1930
$ws.Range("I82").Value = 2040.7693
$ws.Range("K82").Value = 2040.7693
$ws.Range("M82").Value = -1679.7693
$ws.Range("H85").Value = 1930
$ws.Range("I85").Value = 2040.7693
$ws.Range("K85").Value = 2040.7693
$ws.Range("M85").Value = -792.7692999999999

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H109").Value = 69032.75
$ws.Range("I109").Value = 68000
$ws.Range("K109").Value = 68000
$ws.Range("M109").Value = -66613
$ws.Range("H126").Value = 39380572
$ws.Range("I126").Value = 46540090
$ws.Range("K126").Value = 139620270
$ws.Range("M126").Value = -139617800
